# Apply NATMI recompute for Fgf17-Fgfr1: Ligand/Receptor-expressing cell counts
# went from 1 to 3, and all dependent expression/specificity/edge-weight
# columns were recomputed accordingly ("Natmi following Dr Hou advice").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.366566
$ws.Range("H2").Value = 1.099698
$ws.Range("I2").Value = 0.5689653834353526
$ws.Range("J2").Value = 0.5689653834353527
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 5.855348
$ws.Range("N2").Value = 17.566044
$ws.Range("O2").Value = 0.05092948808292105
$ws.Range("P2").Value = 0.05092948808292105
$ws.Range("Q2").Value = 2.146371494968
$ws.Range("R2").Value = 19.317343454712
$ws.Range("S2").Value = 0.0289771157152654
$ws.Range("T2").Value = 0.0289771157152654

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.366566
$ws.Range("H3").Value = 1.099698
$ws.Range("I3").Value = 0.5689653834353526
$ws.Range("J3").Value = 0.5689653834353527
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 62.99699166666667
$ws.Range("N3").Value = 188.990975
$ws.Range("O3").Value = 0.547944295769846
$ws.Range("P3").Value = 0.547944295769846
$ws.Range("Q3").Value = 23.09255524728333
$ws.Range("R3").Value = 207.83299722555
$ws.Range("S3").Value = 0.3117613363439047
$ws.Range("T3").Value = 0.3117613363439048

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.366566
$ws.Range("H4").Value = 1.099698
$ws.Range("I4").Value = 0.5689653834353526
$ws.Range("J4").Value = 0.5689653834353527
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.2640463333333333
$ws.Range("N4").Value = 0.792139
$ws.Range("O4").Value = 0.002296660179179615
$ws.Range("P4").Value = 0.002296660179179615
$ws.Range("Q4").Value = 0.09679040822466667
$ws.Range("R4").Value = 0.8711136740220001
$ws.Range("S4").Value = 0.001306720139467635
$ws.Range("T4").Value = 0.001306720139467635

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.366566
$ws.Range("H5").Value = 1.099698
$ws.Range("I5").Value = 0.5689653834353526
$ws.Range("J5").Value = 0.5689653834353527
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.226170666666667
$ws.Range("N5").Value = 3.678512
$ws.Range("O5").Value = 0.01066516360011862
$ws.Range("P5").Value = 0.01066516360011862
$ws.Range("Q5").Value = 0.4494724765973334
$ws.Range("R5").Value = 4.045252289376
$ws.Range("S5").Value = 0.006068108897142256
$ws.Range("T5").Value = 0.006068108897142256

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.366566
$ws.Range("H6").Value = 1.099698
$ws.Range("I6").Value = 0.5689653834353526
$ws.Range("J6").Value = 0.5689653834353527
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.763133
$ws.Range("N6").Value = 2.289399
$ws.Range("O6").Value = 0.006637687978440185
$ws.Range("P6").Value = 0.006637687978440185
$ws.Range("Q6").Value = 0.279738611278
$ws.Range("R6").Value = 2.517647501502
$ws.Range("S6").Value = 0.003776614685777451
$ws.Range("T6").Value = 0.003776614685777452

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.366566
$ws.Range("H7").Value = 1.099698
$ws.Range("I7").Value = 0.5689653834353526
$ws.Range("J7").Value = 0.5689653834353527
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 43.864011
$ws.Range("N7").Value = 131.592033
$ws.Range("O7").Value = 0.3815267043894945
$ws.Range("P7").Value = 0.3815267043894945
$ws.Range("Q7").Value = 16.079055056226
$ws.Range("R7").Value = 144.711495506034
$ws.Range("S7").Value = 0.2170754876537952
$ws.Range("T7").Value = 0.2170754876537952

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.2777016666666667
$ws.Range("H8").Value = 0.833105
$ws.Range("I8").Value = 0.4310346165646473
$ws.Range("J8").Value = 0.4310346165646473
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 5.855348
$ws.Range("N8").Value = 17.566044
$ws.Range("O8").Value = 0.05092948808292105
$ws.Range("P8").Value = 0.05092948808292105
$ws.Range("Q8").Value = 1.626039898513334
$ws.Range("R8").Value = 14.63435908662
$ws.Range("S8").Value = 0.02195237236765565
$ws.Range("T8").Value = 0.02195237236765565

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.2777016666666667
$ws.Range("H9").Value = 0.833105
$ws.Range("I9").Value = 0.4310346165646473
$ws.Range("J9").Value = 0.4310346165646473
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 62.99699166666667
$ws.Range("N9").Value = 188.990975
$ws.Range("O9").Value = 0.547944295769846
$ws.Range("P9").Value = 0.547944295769846
$ws.Range("Q9").Value = 17.49436958081944
$ws.Range("R9").Value = 157.449326227375
$ws.Range("S9").Value = 0.2361829594259413
$ws.Range("T9").Value = 0.2361829594259413

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.2777016666666667
$ws.Range("H10").Value = 0.833105
$ws.Range("I10").Value = 0.4310346165646473
$ws.Range("J10").Value = 0.4310346165646473
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.2640463333333333
$ws.Range("N10").Value = 0.792139
$ws.Range("O10").Value = 0.002296660179179615
$ws.Range("P10").Value = 0.002296660179179615
$ws.Range("Q10").Value = 0.07332610684388889
$ws.Range("R10").Value = 0.659934961595
$ws.Range("S10").Value = 0.0009899400397119793
$ws.Range("T10").Value = 0.0009899400397119793

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.2777016666666667
$ws.Range("H11").Value = 0.833105
$ws.Range("I11").Value = 0.4310346165646473
$ws.Range("J11").Value = 0.4310346165646473
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 1.226170666666667
$ws.Range("N11").Value = 3.678512
$ws.Range("O11").Value = 0.01066516360011862
$ws.Range("P11").Value = 0.01066516360011862
$ws.Range("Q11").Value = 0.3405096377511111
$ws.Range("R11").Value = 3.06458673976
$ws.Range("S11").Value = 0.004597054702976363
$ws.Range("T11").Value = 0.004597054702976362

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.2777016666666667
$ws.Range("H12").Value = 0.833105
$ws.Range("I12").Value = 0.4310346165646473
$ws.Range("J12").Value = 0.4310346165646473
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.763133
$ws.Range("N12").Value = 2.289399
$ws.Range("O12").Value = 0.006637687978440185
$ws.Range("P12").Value = 0.006637687978440185
$ws.Range("Q12").Value = 0.2119233059883333
$ws.Range("R12").Value = 1.907309753895
$ws.Range("S12").Value = 0.002861073292662734
$ws.Range("T12").Value = 0.002861073292662734

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.2777016666666667
$ws.Range("H13").Value = 0.833105
$ws.Range("I13").Value = 0.4310346165646473
$ws.Range("J13").Value = 0.4310346165646473
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 43.864011
$ws.Range("N13").Value = 131.592033
$ws.Range("O13").Value = 0.3815267043894945
$ws.Range("P13").Value = 0.3815267043894945
$ws.Range("Q13").Value = 12.181108961385
$ws.Range("R13").Value = 109.629980652465
$ws.Range("S13").Value = 0.1644512167356993
$ws.Range("T13").Value = 0.1644512167356993
